# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# Adds a third worksheet "ODI Batting Extra" (sheetId 3) after the existing
# "Player Info" / "ODI Batting" sheets, with extra per-innings batting stats
# for match 4625.

$wb = $excel.ActiveWorkbook

# Duplicate the last existing sheet so the new sheet inherits the same
# sheet-level plumbing (sheetPr/outlinePr, pageMargins, etc.) used
# throughout the workbook, then drop into place right after it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "ODI Batting Extra"
$ws.Cells.Clear()

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws.Cells.Item(1, 2).Value = "BATTING_POSITION"
$ws.Cells.Item(1, 3).Value = "NUM_4"
$ws.Cells.Item(1, 4).Value = "NUM_6"
$ws.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Reuse the same bold/centered/bordered header style already used by the
# other sheets' header rows instead of building a brand-new style.
$headerStyleSource = $wb.Worksheets.Item(1).Range("A1")
$headerStyleSource.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data row (row 2) ---------------------------------------------------
# Everything except BATTING_POSITION is stored as text (even the
# numeric-looking values), so force text entry, then strip the
# now-unneeded "@" number format back off once the literal values are in.
$ws.Range("A2:F2").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "4625"
$ws.Cells.Item(2, 3).Value = "3"
$ws.Cells.Item(2, 4).Value = "0"
$ws.Cells.Item(2, 5).Value = "10.42%"
$ws.Cells.Item(2, 6).Value = "NO"
$ws.Range("A2:F2").ClearFormats()

# BATTING_POSITION is a genuine number.
$ws.Cells.Item(2, 2).Value = 3

# Restore original active sheet/selection.
$wb.Worksheets.Item("Player Info").Activate()
